# Daily IST report: add new submission date column (2026-02-26)
# Mirrors: reports/submissions_daily_matrix.xlsx update (#950)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at L. This pushes the existing L ("total_files")
#     and M ("unique_days") columns one slot to the right (-> M and N),
#     and inherits column L's formatting from the column to its left (K),
#     i.e. the plain bold/centered date-header style.
$ws.Columns("L").Insert()

# --- Column widths: L becomes a normal date-data column (width 12),
#     M keeps the "total_files" header width (13), N (new) gets the
#     "unique_days" header width (13). ColumnWidth (Excel chars) is the
#     stored xlsx <col width> minus a fixed 0.83 offset.
$ws.Range("L1").ColumnWidth = 11.17
$ws.Range("M1").ColumnWidth = 12.17
$ws.Range("N1").ColumnWidth = 12.17

# --- Header row: L1 becomes the new date "2026-02-26" (same text style as
#     the other date headers). Writing the literal string directly makes
#     Excel auto-detect it as a date and reformat the cell, so instead we
#     write it as a text formula and then paste-special as a value, which
#     keeps it a plain text cell using the inherited date-header style.
$ws.Range("L1").Formula = '="2026-02-26"'
$ws.Range("L1").Copy()
$ws.Range("L1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# M1/N1 already hold the shifted "total_files"/"unique_days" labels after
# the column insert, so nothing else to do for the header row.

# --- Data rows 2-109: new submissions for 2026-02-26 (1 = submitted a
#     file that day, 0 = no submission). total_files (M) and unique_days
#     (N) are recomputed from their pre-shift values (which the insert
#     already carried into M/N) plus this new day's contribution.
$newDay = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,1,1,0,1,1,0,0,0,1,1,1,1,1,1,1,0,0,1,1,1,0,1,1,1,1,0,1,1,1,0,1,0,1,0,1,1,0,1,0,0,1,1,0,0,0,0,0,1,0,0,1,1,1,1,0,0,1,1,0,1,1,0,0,1,0,0,0,1,0,0,1,1,1,0,0,1,0,0,1,0,1,0,1,1,0,1,0)

for ($i = 0; $i -lt $newDay.Length; $i++) {
    $r = $i + 2
    $flag = $newDay[$i]

    $ws.Cells.Item($r, 12).Value2 = $flag          # L: 2026-02-26 flag

    $priorTotal = $ws.Cells.Item($r, 13).Value()   # M currently holds old total_files
    $priorDays  = $ws.Cells.Item($r, 14).Value()   # N currently holds old unique_days

    $ws.Cells.Item($r, 13).Value2 = $priorTotal + $flag
    $ws.Cells.Item($r, 14).Value2 = $priorDays + ($(if ($flag -gt 0) { 1 } else { 0 }))
}
